$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.641.34'
$ws.Range('E2').Value = '  +1.24%  '
$ws.Range('D3').Value = '1.826.60'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.44%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.007'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '308.67'
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('E7').Value = '  +3.90%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3607'
$ws.Range('E8').Value = '  +0.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07133'
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9031'
$ws.Range('E10').Value = '  +2.38%  '
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').Value = '1.852.03'
$ws.Range('E13').Value = '  +2.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.276'
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.361'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.54'
$ws.Range('E16').Value = '  +3.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.009'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008550'
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.008'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').Value = '26.676.88'
$ws.Range('E20').Value = '  +1.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.22'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.025'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('E24').Value = '  -3.33%  '
$ws.Range('E25').Value = '  +1.29%  '
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.979'
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '113.93'
$ws.Range('E28').Value = '  +1.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.869'
$ws.Range('E29').Value = '  -0.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08812'
$ws.Range('E30').Value = '  +1.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.146'
$ws.Range('E31').Value = '  +2.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.846'
$ws.Range('E32').Value = '  +5.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.167'
$ws.Range('E33').Value = '  +5.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7371'
$ws.Range('E34').Value = '  +2.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.445'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01932'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05161'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.881'
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5057'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.046'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4668'
$ws.Range('E45').Value = '  +1.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.980'
$ws.Range('E46').Value = '  +0.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '98.02'
$ws.Range('E47').Value = '  -3.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.574'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06049'
$ws.Range('E49').Value = '  +1.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '64.04'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.80'
$ws.Range('E51').Value = '  -0.22%  '
